$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap row 3 (Jalen Green / PG,SG / Houston Rockets) with row 14 (Klay Thompson / SG,SF / Dallas Mavericks)
$a3 = $ws.Cells.Item(3,1).Value2
$b3 = $ws.Cells.Item(3,2).Value2
$c3 = $ws.Cells.Item(3,3).Value2

$ws.Cells.Item(3,1).Value = $ws.Cells.Item(14,1).Value2
$ws.Cells.Item(3,2).Value = $ws.Cells.Item(14,2).Value2
$ws.Cells.Item(3,3).Value = $ws.Cells.Item(14,3).Value2

$ws.Cells.Item(14,1).Value = $a3
$ws.Cells.Item(14,2).Value = $b3
$ws.Cells.Item(14,3).Value = $c3

# Swap row 6 (Pascal Siakam / SF,PF / Indiana Pacers) with row 15 (Jaylen Brown / SG,SF / Boston Celtics)
$a6 = $ws.Cells.Item(6,1).Value2
$b6 = $ws.Cells.Item(6,2).Value2
$c6 = $ws.Cells.Item(6,3).Value2

$ws.Cells.Item(6,1).Value = $ws.Cells.Item(15,1).Value2
$ws.Cells.Item(6,2).Value = $ws.Cells.Item(15,2).Value2
$ws.Cells.Item(6,3).Value = $ws.Cells.Item(15,3).Value2

$ws.Cells.Item(15,1).Value = $a6
$ws.Cells.Item(15,2).Value = $b6
$ws.Cells.Item(15,3).Value = $c6
